# Commit: "Mult new tests added (for pose extraction verification)"
#
# Target sheet layout (tab order):
#   Main, roll_and_pitch_new, Lin_New, Rotation_new, Test_1_rotation,
#   Test_1a, Test_1b, Test_1c
#
# The original "Test_1a", "Test_1b", "Test_1c" sheets are duplicated to the
# end of the workbook (keeping their original content) to become the new
# "Test_1a" / "Test_1b" / "Test_1c" tabs. The original four sheets
# ("Test_1_rotation", "Test_1a", "Test_1b", "Test_1c") are then renamed and
# their Parameter/Info cells updated in place to become the four new test
# configs ("roll_and_pitch_new", "Lin_New", "Rotation_new",
# "Test_1_rotation"), after which everything is reordered.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate Test_1a / Test_1b / Test_1c to the end of the workbook
# BEFORE touching their contents, so the copies retain the original data.
# These duplicates will become the "new" Test_1a / Test_1b / Test_1c.
# ---------------------------------------------------------------------
$srcA = $wb.Worksheets.Item("Test_1a")
$srcA.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dupA = $wb.Worksheets.Item($wb.Worksheets.Count)

$srcB = $wb.Worksheets.Item("Test_1b")
$srcB.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dupB = $wb.Worksheets.Item($wb.Worksheets.Count)

$srcC = $wb.Worksheets.Item("Test_1c")
$srcC.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dupC = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------
# Step 2: turn the original "Test_1_rotation" sheet into "roll_and_pitch_new"
# ---------------------------------------------------------------------
$rp = $wb.Worksheets.Item("Test_1_rotation")
$rp.Name = "roll_and_pitch_new"
$rp.Range("B2").Value = "roll_and_pitch_new.txt"
$rp.Range("B6").Value = "KahnPhone_new.sdf"
$rp.Range("B7").Value = "DICT_4X4_50_s100_id0.sdf"
$rp.Range("B8").Value = "DICT_4X4_50_s100_id1.sdf"
$rp.Range("C8").Value = "-1,-1,0,0,0,0"

# ---------------------------------------------------------------------
# Step 3: turn the original "Test_1a" sheet into "Lin_New"
# ---------------------------------------------------------------------
$lin = $wb.Worksheets.Item("Test_1a")
$lin.Name = "Lin_New"
$lin.Range("B2").Value = "lin_move_new.txt"
$lin.Range("B6").Value = "KahnPhone_new.sdf"
$lin.Range("B7").Value = "DICT_4X4_50_s100_id0.sdf"
$lin.Range("B8").Value = "DICT_4X4_50_s100_id1.sdf"
$lin.Range("C8").Value = "-1,-1,0,0,0,0"

# ---------------------------------------------------------------------
# Step 4: turn the original "Test_1b" sheet into "Rotation_new"
# ---------------------------------------------------------------------
$rot = $wb.Worksheets.Item("Test_1b")
$rot.Name = "Rotation_new"
$rot.Range("B2").Value = "rotate_new.txt"
$rot.Range("B6").Value = "KahnPhone_new.sdf"
$rot.Range("B7").Value = "DICT_4X4_50_s100_id0.sdf"
$rot.Range("B8").Value = "DICT_4X4_50_s100_id1.sdf"
$rot.Range("C8").Value = "-1,-1,0,0,0,0"

# ---------------------------------------------------------------------
# Step 5: turn the original "Test_1c" sheet into the new "Test_1_rotation"
# (restores the values the old "Test_1_rotation" sheet used to have)
# ---------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Test_1c")
$tr.Name = "Test_1_rotation"
$tr.Range("B2").Value = "rotate.txt"
$tr.Range("B6").Value = "KahnPhone.sdf"
$tr.Range("B7").Value = "DICT_4X4_50_s100_id0.sdf"
$tr.Range("B8").Value = "DICT_4X4_50_s100_id1.sdf"
$tr.Range("C8").Value = "-1,-1,0,0,0,0"

# ---------------------------------------------------------------------
# Step 6: rename the duplicated sheets back to Test_1a / Test_1b / Test_1c
# ---------------------------------------------------------------------
$dupA.Name = "Test_1a"
$dupB.Name = "Test_1b"
$dupC.Name = "Test_1c"

# ---------------------------------------------------------------------
# Step 7: re-order all the sheets into the final tab order
# ---------------------------------------------------------------------
$main = $wb.Worksheets.Item("Main")

$rpWs = $wb.Worksheets.Item("roll_and_pitch_new")
$rpWs.Move($null, $main)

$linWs = $wb.Worksheets.Item("Lin_New")
$linWs.Move($null, $rpWs)

$rotWs = $wb.Worksheets.Item("Rotation_new")
$rotWs.Move($null, $linWs)

$trWs = $wb.Worksheets.Item("Test_1_rotation")
$trWs.Move($null, $rotWs)

$aWs = $wb.Worksheets.Item("Test_1a")
$aWs.Move($null, $trWs)

$bWs = $wb.Worksheets.Item("Test_1b")
$bWs.Move($null, $aWs)

$cWs = $wb.Worksheets.Item("Test_1c")
$cWs.Move($null, $bWs)

foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
